$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(256, 44330, 1, 4, 64.69351447517387)
    ,@(257, 44331, 0, 4, 64.69351447517387)
    ,@(258, 44332, 1, 4, 64.69351447517387)
    ,@(259, 44333, 0, 3, 48.5201358563804)
    ,@(260, 44334, 0, 3, 48.5201358563804)
    ,@(261, 44335, 0, 2, 32.34675723758694)
    ,@(262, 44336, 2, 4, 64.69351447517387)
    ,@(263, 44337, 0, 3, 48.5201358563804)
    ,@(264, 44338, 0, 3, 48.5201358563804)
    ,@(265, 44339, 0, 2, 32.34675723758694)
    ,@(266, 44340, 3, 5, 80.86689309396733)
    ,@(267, 44341, 0, 5, 80.86689309396733)
    ,@(268, 44342, 0, 5, 80.86689309396733)
    ,@(269, 44343, 0, 3, 48.5201358563804)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Range("A255").Copy($ws.Range("A$r"))
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
}
